$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$data = @{
    2  = @{ D = 44320; M = 80;  N = 16000; O = 17000; P = 16500; S = 825 }
    3  = @{ D = 44792; M = 100; N = 21000; O = 22000; P = 21500; S = 1075 }
    4  = @{ D = 44533; M = 100; N = 16000; O = 17000; P = 16500; S = 825 }
    5  = @{ D = 44761; M = 100; N = 20000; O = 21000; P = 20500; S = 1025 }
    6  = @{ D = 44890; M = 80;  N = 20000; O = 23000; P = 22250; S = 1112 }
    7  = @{ D = 44357; M = 100; N = 14000; O = 15000; P = 14500; S = 725 }
    8  = @{ D = 44893; M = 80;  N = 21000; O = 22000; P = 21625; S = 1081 }
    9  = @{ D = 45092; M = 150; N = 24000; O = 25000; P = 24333; S = 1217 }
    10 = @{ D = 44798; M = 80;  N = 21000; O = 22000; P = 21500; S = 1075 }
    11 = @{ D = 44708; M = 80;  N = 20000; O = 21000; P = 20500; S = 1025 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
